# Remove the trailing site-footer block that was scraped along with the
# bibliography ("Ver no Jupiter Salvar em pdf Salvar em docx" plus the
# "© 2020 ... Creative Commons Attribution" copyright line), together with
# the blank paragraph that separated it from the bibliography text. The
# blank paragraph / page-break paragraph that come after the footer are
# left untouched.

$d = $word.ActiveDocument

$footer = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $footer = $p
        break
    }
}

if ($footer -ne $null) {
    $copyright = $footer.Next()
    $blankBefore = $footer.Previous()

    $startPos = $blankBefore.Range.Start
    $endPos = $copyright.Range.End

    $r = $d.Range($startPos, $endPos)
    $r.Delete()
}
